# untypedDecimals.xlsx fix for failing atomic tabOTTR tests.
#
# The "data"/1/"auto" triple in column A (rows 7-9) shifts down by one:
#   A7 "data" (default style)      -> A7 becomes the number 1 (keeps default style)
#   A8 the number 1 (styled s=1)   -> A8 becomes "auto" (keeps styled s=1)
#   A9 "auto" (styled s=1)         -> A9 becomes "data" (reverts to default style)
# The active selection follows the content up from A9 to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: becomes "data" and drops back to the default (unstyled) look that
# A7 originally had -- explicitly re-assert the workbook's base font so the
# engine reuses the existing default style instead of minting a new one.
$ws.Range("A9").Value = "data"
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 10

# Row 8: becomes "auto", keeping its existing (styled) formatting untouched.
$ws.Range("A8").Value = "auto"

# Row 7: becomes the numeric value 1, keeping its existing default formatting.
$ws.Range("A7").Value = 1

# The sheet's active/selected cell moves from A9 to A8.
$ws.Range("A8").Select() | Out-Null

# Cosmetic: bump the tab-bar/scrollbar split ratio slightly (980 -> 990).
$excel.ActiveWindow.TabRatio = 0.99
